# daily auto push: 2026-01-12 09:44 UTC
# Insert a new data row for 2026/01/12 (time 16) just before the
# 2026/12/29 block, keeping the sheet's chronological row order.
# All rows from the old row 609 ("2026/12/29 ...") down to the end
# shift down by one row, and the new row's data lands at row 609.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything at/after row 609 down by one row.
$ws.Rows.Item(609).Insert()

# Populate the newly-opened row 609. The date/day-of-week columns are
# stored as plain text in this sheet (e.g. "2026/01/12"), so format the
# date cell as text before assigning it to stop Excel from silently
# re-typing it as a date serial, then reset the style back to Normal so
# no stray number-format style id is left behind on the cell.
$ws.Cells.Item(609, 1).NumberFormat = "@"
$ws.Cells.Item(609, 1).Value = "2026/01/12"
$ws.Cells.Item(609, 1).Style = "Normal"

$ws.Cells.Item(609, 2).Value = "月"
$ws.Cells.Item(609, 3).Value = 16
$ws.Cells.Item(609, 4).Value = 24
